$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Update the "Last Updated" DATE field result text in the footer.
# ---------------------------------------------------------------------
$footer = $d.Sections(1).Footers(1)
$footer.Range.Find.Execute("20/02/2023 16:13", $true, $false, $false, $false, $false, $true, 1, $false, "25/02/2023 16:29", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Append the new "Version 0.0.4" changelog section at the end of the
#    document body (after the existing "...2000Mb" bullet, before the
#    section properties).
# ---------------------------------------------------------------------

# Grab a reference to the list-numbering template already used by the
# other "ListParagraph" bullets (numId 5) so new bullets continue it
# rather than starting a brand-new list.
$bulletTemplate = $d.Paragraphs($d.Paragraphs.Count).Range.ListFormat.ListTemplate

function New-TrailingParagraph {
    # Inserts a new, *truly* empty paragraph (no residual run) after the
    # current last paragraph and returns it. Plain InsertParagraphAfter()
    # leaves a stray empty <w:r/> behind, so instead a throw-away
    # character is typed and then deleted again, which leaves a clean
    # paragraph mark with no run at all.
    $last = $d.Paragraphs($d.Paragraphs.Count)
    $r = $last.Range
    $r.Collapse(0)
    $r.InsertAfter("`rx")
    $newPara = $d.Paragraphs($d.Paragraphs.Count)
    $junk = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)
    $junk.Delete()
    return $d.Paragraphs($d.Paragraphs.Count)
}

function Add-RunText {
    # Appends $text as its own run at the end of the document, using a
    # throw-away bookmark to stop it from being coalesced into the
    # previous run.
    param([string]$text)
    $r = $d.Paragraphs($d.Paragraphs.Count).Range
    $r.Collapse(0)
    $mark = $d.Bookmarks.Add("zzTmpRunBreak", $r)
    $r.InsertAfter($text)
    $d.Bookmarks("zzTmpRunBreak").Delete()
}

# --- blank paragraph -----------------------------------------------
$p = New-TrailingParagraph
$p.Range.Style = $d.Styles("Normal")

# --- Heading 1: "Version 0.0.4" -------------------------------------
$p = New-TrailingParagraph
$p.Range.Style = $d.Styles("Heading 1")
Add-RunText "Version 0.0."
Add-RunText "4"

# --- Heading 2: "Added" ----------------------------------------------
$p = New-TrailingParagraph
$p.Range.Style = $d.Styles("Heading 2")
Add-RunText "Added"

# --- Heading 3: "Items" ----------------------------------------------
$p = New-TrailingParagraph
$p.Range.Style = $d.Styles("Heading 3")
Add-RunText "Items"

# --- Heading 2: "Fixed" ------------------------------------------------
$p = New-TrailingParagraph
$p.Range.Style = $d.Styles("Heading 2")
Add-RunText "Fixed"

# --- Heading 2: "Updated" ----------------------------------------------
$p = New-TrailingParagraph
$p.Range.Style = $d.Styles("Heading 2")
Add-RunText "Updated"

# --- ListParagraph bullet: "Texture Streaming Size to 4096Mb" ----------
$p = New-TrailingParagraph
$p.Range.Style = $d.Styles("List Paragraph")
$p.Range.ListFormat.ApplyListTemplate($bulletTemplate, $true)
Add-RunText "Texture Streaming Size to "
Add-RunText "4"
Add-RunText "0"
Add-RunText "96"
Add-RunText "Mb"

# --- trailing blank paragraph -------------------------------------------
$p = New-TrailingParagraph
$p.Range.Style = $d.Styles("Normal")

Write-Output "done"
